# Fruta / hortaliza, semanal
# Insert a new weekly record row just above the current row 92 (Espárragos,
# Macroferia Regional de Talca), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 92 - this shifts rows 92..118 down to 93..119
# and keeps formatting (e.g. the date style on column D) consistent with the
# row that used to occupy that position.
$ws.Range("A92").EntireRow.Insert()

# Populate the newly inserted row 92 with the new weekly record.
$ws.Cells.Item(92, 1).Value = 5
$ws.Cells.Item(92, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(92, 3).Value = "Maule"
$ws.Cells.Item(92, 4).Value = 45218
$ws.Cells.Item(92, 5).Value = 7
$ws.Cells.Item(92, 6).Value = 300000000
$ws.Cells.Item(92, 7).Value = "Espárragos"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 3000
$ws.Cells.Item(92, 11).Value = 1200
$ws.Cells.Item(92, 12).Value = 1200
$ws.Cells.Item(92, 13).Value = 1200
$ws.Cells.Item(92, 14).Value = "`$/kilo"
$ws.Cells.Item(92, 15).Value = "Provincia de Linares"
$ws.Cells.Item(92, 16).Value = 1200
$ws.Cells.Item(92, 17).Value = 1
$ws.Cells.Item(92, 18).Value = "Hortaliza"
